# Reorder header labels in row 1 (columns C:F) on every worksheet.
# Old order: C=normalize_group, D=trajgroup_no_vary_q, E=uniform_scaling_q, F=variable_trajectory_group
# New order: C=variable_trajectory_group, D=normalize_group, E=trajgroup_no_vary_q, F=uniform_scaling_q
# (i.e. "variable_trajectory_group" moves from F1 to C1, pushing the other three one column right)

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("C1").Value = "variable_trajectory_group"
    $ws.Range("D1").Value = "normalize_group"
    $ws.Range("E1").Value = "trajgroup_no_vary_q"
    $ws.Range("F1").Value = "uniform_scaling_q"
}
